# Add 2022-Q4 data
# 1. Insert a new worksheet named "2022-Q4" right before the "2022-Q3" sheet
#    and populate it with the quarterly fund-holding table.
# 2. Insert a new row into the "总计" (summary) sheet with the 2022-Q4 totals,
#    pushing the existing quarters down by one row.
#
# NOTE: sheet object variables in this COM layer are positional handles, not
# stable identities -- once the sheet collection is reshuffled by Add()/
# Insert(), any previously-fetched $sheet variable can silently start
# pointing at a different sheet. So: do ALL structural changes (adding /
# inserting / renaming sheets & rows) first, and only AFTER that is finished
# re-fetch the sheet references we need and perform the data writes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1 (structural): create the "2022-Q4" sheet positioned before "2022-Q3"
# ---------------------------------------------------------------------------
$q3Anchor = $wb.Sheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($q3Anchor)
$newSheet.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# Step 2 (structural): insert the 2022-Q4 row into the "总计" summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Sheets.Item("总计")
$summary.Rows.Item(2).Insert()

# ---------------------------------------------------------------------------
# Everything below is pure data/formatting -- re-fetch fresh references.
# ---------------------------------------------------------------------------
$q4sheet = $wb.Sheets.Item("2022-Q4")
$q3sheet = $wb.Sheets.Item("2022-Q3")
# "2022-Q2" already has exactly 9 data rows (A2:A10) styled the same way the
# new sheet needs, so borrow its formatting for the bold index column.
$styleDonor = $wb.Sheets.Item("2022-Q2")
$summary2 = $wb.Sheets.Item("总计")

# --- 2022-Q4 sheet: header row -----------------------------------------
$q4sheet.Cells.Item(1, 2).Value = "基金代码"
$q4sheet.Cells.Item(1, 3).Value = "基金名称"
$q4sheet.Cells.Item(1, 4).Value = "基金规模"
$q4sheet.Cells.Item(1, 5).Value = "股票总仓位"
$q4sheet.Cells.Item(1, 6).Value = "仓位占比"
$q4sheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4sheet.Cells.Item(1, 8).Value = "仓位排名"

# --- 2022-Q4 sheet: data rows -------------------------------------------
# A=index(number), B=基金代码(text), C=基金名称(text), D=基金规模(text),
# E=股票总仓位(text), F=仓位占比(text), G=持有市值(text), H=仓位排名(number)
$rows = @(
    @(0, "001643", "汇丰晋信智造先锋股票A", "16.39", "93.66", "6.95", "1.1391", 1),
    @(1, "011578", "汇丰晋信核心成长混合A", "21.58", "94.18", "4.19", "0.9042", 9),
    @(2, "001644", "汇丰晋信智造先锋股票C", "8.77", "93.66", "6.95", "0.6095", 1),
    @(3, "011579", "汇丰晋信核心成长混合C", "4.18", "94.18", "4.19", "0.1751", 9),
    @(4, "217021", "招商优势企业混合", "3.40", "79.82", "4.63", "0.1574", 9),
    @(5, "014320", "德邦半导体产业混合C", "1.52", "92.57", "6.85", "0.1041", 2),
    @(6, "015071", "鑫元专精特新混合A", "2.46", "70.69", "2.62", "0.0645", 10),
    @(7, "014319", "德邦半导体产业混合A", "0.37", "92.57", "6.85", "0.0253", 2),
    @(8, "015072", "鑫元专精特新混合C", "0.25", "70.69", "2.62", "0.0066", 10)
)

$r = 2
foreach ($row in $rows) {
    $q4sheet.Cells.Item($r, 1).Value = $row[0]
    $q4sheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $q4sheet.Cells.Item($r, 3).Value = $row[2]
    $q4sheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $q4sheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $q4sheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $q4sheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $q4sheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# --- 2022-Q4 sheet: formatting to match the other quarter sheets -------
# Bold/centred header (B1:H1).
for ($col = 2; $col -le 8; $col++) {
    $q3sheet.Cells.Item(1, $col).Copy()
    $q4sheet.Cells.Item(1, $col).PasteSpecial(-4122)
}
# Bold index column (A2:A10) -- "2022-Q2" has the same 9-row shape.
for ($rr = 2; $rr -le 10; $rr++) {
    $styleDonor.Cells.Item($rr, 1).Copy()
    $q4sheet.Cells.Item($rr, 1).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 总计 sheet: fill in the inserted row with the 2022-Q4 totals
# ---------------------------------------------------------------------------
$summary2.Cells.Item(2, 1).Value = 0
$summary2.Cells.Item(2, 2).Value = "2022-Q4"
$summary2.Cells.Item(2, 3).Value = 9
$summary2.Cells.Item(2, 4).Value = 3.19

# The inserted row's B:D cells pick up stray formatting from the insert;
# clear it back to the default (unstyled) look used by the other data rows.
$summary2.Range("B2:D2").Style = "Normal"

# Restore the bold index-column styling on A2 (Insert() drops it).
$summary2.Cells.Item(3, 1).Copy()
$summary2.Cells.Item(2, 1).PasteSpecial(-4122)
$summary2.Cells.Item(2, 1).Value = 0
